$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "last updated" timestamp (row 1)
$ws.Range("A1").Value = "Datos actualizados a 17 de Junio de 2020 a las 08:33"

# --- Country name (column A) reorderings ---
$ws.Range("A41").Value = "Afganistan"
$ws.Range("A42").Value = "Filipinas"

$ws.Range("A50").Value = "Israel"
$ws.Range("A51").Value = "Barein"

$ws.Range("A130").Value = "Georgia"
$ws.Range("A131").Value = "Yemen"
$ws.Range("A132").Value = "Congo"

$ws.Range("A206").Value = "Islas Malvinas"
$ws.Range("A207").Value = "Groenlandia"

$ws.Range("A210").Value = "Montserrat"
$ws.Range("A211").Value = "Seychelles"

$ws.Range("A213").Value = "Islas Virgenes Britanicas"
$ws.Range("A214").Value = "Papua Nueva Guinea"

# --- Updated statistics (columns B-H) ---

# Row 13: Alemania
$ws.Range("D13").Value = 173600
$ws.Range("E13").Value = 5872

# Row 38: Ucrania
$ws.Range("B38").Value = 33234
$ws.Range("C38").Value = 758
$ws.Range("D38").Value = 14943
$ws.Range("E38").Value = 17348
$ws.Range("G38").Value = 31
$ws.Range("H38").Value = 943

# Row 41: Afganistan (new values)
$ws.Range("B41").Value = 26874
$ws.Range("C41").Value = 564
$ws.Range("D41").Value = 6158
$ws.Range("E41").Value = 20212
$ws.Range("G41").Value = 13
$ws.Range("H41").Value = 504

# Row 42: Filipinas (shifted values)
$ws.Range("B42").Value = 26781
$ws.Range("D42").Value = 6552
$ws.Range("E42").Value = 19126
$ws.Range("H42").Value = 1103

# Row 50: Israel (shifted values)
$ws.Range("B50").Value = 19637
$ws.Range("C50").Value = 142
$ws.Range("D50").Value = 15459
$ws.Range("E50").Value = 3875
$ws.Range("G50").Value = 1
$ws.Range("H50").Value = 303

# Row 51: Barein (shifted values)
$ws.Range("B51").Value = 19553
$ws.Range("D51").Value = 13866
$ws.Range("E51").Value = 5640
$ws.Range("H51").Value = 47

# Row 129: Burkina Faso
$ws.Range("D129").Value = 809
$ws.Range("E129").Value = 33

# Row 130: Georgia (new values)
$ws.Range("B130").Value = 888
$ws.Range("C130").Value = 9
$ws.Range("D130").Value = 731
$ws.Range("E130").Value = 143
$ws.Range("H130").Value = 14

# Row 131: Yemen (shifted values)
$ws.Range("B131").Value = 885
$ws.Range("D131").Value = 91
$ws.Range("E131").Value = 580
$ws.Range("H131").Value = 214

# Row 132: Congo (shifted values)
$ws.Range("B132").Value = 883
$ws.Range("D132").Value = 391
$ws.Range("E132").Value = 465
$ws.Range("H132").Value = 27

# Row 210: Montserrat (shifted values)
$ws.Range("D210").Value = 10
$ws.Range("H210").Value = 1

# Row 211: Seychelles (shifted values)
$ws.Range("D211").Value = 11
$ws.Range("H211").Value = 0

# Row 213: Islas Virgenes Britanicas (shifted values)
$ws.Range("D213").Value = 7
$ws.Range("H213").Value = 1

# Row 214: Papua Nueva Guinea (shifted values)
$ws.Range("D214").Value = 8
$ws.Range("H214").Value = 0
